$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 18:22"

# --- Etiopia overtakes Dinamarca in the ranking (rows 77-78 swap labels) ---
$ws.Cells.Item(77, 1).Value = "Etiopia"
$ws.Cells.Item(78, 1).Value = "Dinamarca"

# --- Refresh per-country COVID figures (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Cells.Item(4, 2).Value = 4336922
$ws.Cells.Item(4, 3).Value = 21213
$ws.Cells.Item(4, 4).Value = 2063411
$ws.Cells.Item(4, 5).Value = 2123954
$ws.Cells.Item(4, 7).Value = 159
$ws.Cells.Item(4, 8).Value = 149557
$ws.Cells.Item(5, 2).Value = 2402255
$ws.Cells.Item(5, 3).Value = 5821
$ws.Cells.Item(5, 5).Value = 698184
$ws.Cells.Item(5, 7).Value = 95
$ws.Cells.Item(5, 8).Value = 86591
$ws.Cells.Item(6, 2).Value = 1428229
$ws.Cells.Item(6, 3).Value = 42735
$ws.Cells.Item(6, 4).Value = 913657
$ws.Cells.Item(6, 5).Value = 481849
$ws.Cells.Item(6, 7).Value = 627
$ws.Cells.Item(6, 8).Value = 32723
$ws.Cells.Item(13, 2).Value = 299426
$ws.Cells.Item(13, 3).Value = 745
$ws.Cells.Item(13, 7).Value = 14
$ws.Cells.Item(13, 8).Value = 45752
$ws.Cells.Item(17, 2).Value = 246118
$ws.Cells.Item(17, 3).Value = 254
$ws.Cells.Item(17, 4).Value = 198446
$ws.Cells.Item(17, 5).Value = 12565
$ws.Cells.Item(17, 7).Value = 5
$ws.Cells.Item(17, 8).Value = 35107
$ws.Cells.Item(24, 2).Value = 113862
$ws.Cells.Item(24, 3).Value = 306
$ws.Cells.Item(24, 4).Value = 99344
$ws.Cells.Item(24, 5).Value = 5628
$ws.Cells.Item(24, 7).Value = 5
$ws.Cells.Item(24, 8).Value = 8890
$ws.Cells.Item(44, 2).Value = 52946
$ws.Cells.Item(44, 3).Value = 214
$ws.Cells.Item(71, 2).Value = 15273
$ws.Cells.Item(71, 3).Value = 61
$ws.Cells.Item(71, 4).Value = 11423
$ws.Cells.Item(71, 5).Value = 3479
$ws.Cells.Item(71, 7).Value = 2
$ws.Cells.Item(71, 8).Value = 371
$ws.Cells.Item(77, 2).Value = 13968
$ws.Cells.Item(77, 3).Value = 720
$ws.Cells.Item(77, 4).Value = 6216
$ws.Cells.Item(77, 5).Value = 7529
$ws.Cells.Item(77, 7).Value = 14
$ws.Cells.Item(77, 8).Value = 223
$ws.Cells.Item(78, 2).Value = 13438
$ws.Cells.Item(78, 4).Value = 12340
$ws.Cells.Item(78, 5).Value = 485
$ws.Cells.Item(78, 8).Value = 613
$ws.Cells.Item(92, 2).Value = 7192
$ws.Cells.Item(92, 3).Value = 42
$ws.Cells.Item(92, 4).Value = 5970
$ws.Cells.Item(92, 5).Value = 1163
$ws.Cells.Item(99, 2).Value = 4763
$ws.Cells.Item(99, 3).Value = 126
$ws.Cells.Item(99, 4).Value = 2682
$ws.Cells.Item(99, 5).Value = 1943
$ws.Cells.Item(99, 7).Value = 4
$ws.Cells.Item(99, 8).Value = 138
$ws.Cells.Item(104, 2).Value = 4193
$ws.Cells.Item(104, 3).Value = 27
$ws.Cells.Item(104, 5).Value = 2617
$ws.Cells.Item(104, 7).Value = 1
$ws.Cells.Item(104, 8).Value = 202
$ws.Cells.Item(114, 2).Value = 2782
$ws.Cells.Item(114, 3).Value = 12
$ws.Cells.Item(114, 5).Value = 665
$ws.Cells.Item(135, 2).Value = 1669
$ws.Cells.Item(135, 3).Value = 53
$ws.Cells.Item(135, 4).Value = 593
$ws.Cells.Item(135, 5).Value = 1065
$ws.Cells.Item(145, 5).Value = 119
$ws.Cells.Item(145, 7).Value = 1
$ws.Cells.Item(145, 8).Value = 2
